$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 14370.454
$ws.Range("J87").Value = 14370.454
$ws.Range("L87").Value = 14370.454
$ws.Range("N87").Value = -16866.454

$ws.Range("H90").Value = 14370.454
$ws.Range("J90").Value = 14370.454
$ws.Range("L90").Value = 43111.362
$ws.Range("N90").Value = -55591.362

$ws.Range("H107").Value = 339.96667
$ws.Range("I107").Value = 301
$ws.Range("J107").Value = 534.8
$ws.Range("K107").Value = 301
$ws.Range("L107").Value = 534.8
$ws.Range("M107").Value = 1619
$ws.Range("N107").Value = -4374.8

$ws.Range("H125").Value = 7659.1333
$ws.Range("I125").Value = 9753.362999999999
$ws.Range("K125").Value = 87780.26699999999
$ws.Range("M125").Value = -85320.26699999999

$ws.Range("H131").Value = 7497.0527
$ws.Range("I131").Value = 12581.556
$ws.Range("J131").Value = 2921
$ws.Range("K131").Value = 37744.66800000001
$ws.Range("L131").Value = 8763
$ws.Range("M131").Value = -32704.66800000001
$ws.Range("N131").Value = -18843

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3583.2703
$ws.Range("I32").Value = 2566.6765
$ws.Range("J32").Value = 15104.667
$ws.Range("K32").Value = 2566.6765
$ws.Range("L32").Value = 15104.667
$ws.Range("M32").Value = -2279.6765
$ws.Range("N32").Value = -15678.667

$ws.Range("H45").Value = 1961.7273
$ws.Range("I45").Value = 2007.9
$ws.Range("K45").Value = 2007.9
$ws.Range("M45").Value = -1630.9

$ws.Range("H88").Value = 2298
$ws.Range("I88").Value = 2400
$ws.Range("J88").Value = 2272.5
$ws.Range("K88").Value = 2400
$ws.Range("L88").Value = 2272.5
$ws.Range("M88").Value = -1994
$ws.Range("N88").Value = -3084.5

$ws.Range("H91").Value = 2298
$ws.Range("I91").Value = 2400
$ws.Range("J91").Value = 2272.5
$ws.Range("K91").Value = 2400
$ws.Range("L91").Value = 2272.5
$ws.Range("M91").Value = -996
$ws.Range("N91").Value = -5080.5

$ws.Range("H97").Value = 1160.875
$ws.Range("I97").Value = 934.7037
$ws.Range("J97").Value = 2382.2
$ws.Range("K97").Value = 934.7037
$ws.Range("L97").Value = 2382.2
$ws.Range("M97").Value = -438.7037
$ws.Range("N97").Value = -3374.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 430062.72
$ws.Range("I64").Value = 937586.6
$ws.Range("J64").Value = 619.38464
$ws.Range("K64").Value = 937586.6
$ws.Range("L64").Value = 619.38464
$ws.Range("M64").Value = -937361.6
$ws.Range("N64").Value = -1069.38464

$ws.Range("H67").Value = 430062.72
$ws.Range("I67").Value = 937586.6
$ws.Range("J67").Value = 619.38464
$ws.Range("K67").Value = 937586.6
$ws.Range("L67").Value = 619.38464
$ws.Range("M67").Value = -936806.6
$ws.Range("N67").Value = -2179.38464

$ws.Range("H86").Value = 3402.6428
$ws.Range("I86").Value = 3419.1667
$ws.Range("J86").Value = 3303.5
$ws.Range("K86").Value = 3419.1667
$ws.Range("L86").Value = 3303.5
$ws.Range("M86").Value = -2296.1667
$ws.Range("N86").Value = -5549.5

$ws.Range("H89").Value = 3402.6428
$ws.Range("I89").Value = 3419.1667
$ws.Range("J89").Value = 3303.5
$ws.Range("K89").Value = 17095.8335
$ws.Range("L89").Value = 16517.5
$ws.Range("M89").Value = -11479.8335
$ws.Range("N89").Value = -27749.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 5672
$ws.Range("I62").Value = 3815.4285
$ws.Range("K62").Value = 3815.4285
$ws.Range("M62").Value = -3191.4285

$ws.Range("H65").Value = 5672
$ws.Range("I65").Value = 3815.4285
$ws.Range("K65").Value = 19077.1425
$ws.Range("M65").Value = -15957.1425

$ws.Range("H132").Value = 3280.72
$ws.Range("I132").Value = 3672.9167
$ws.Range("J132").Value = 2272.2144
$ws.Range("K132").Value = 11018.7501
$ws.Range("L132").Value = 6816.6432
$ws.Range("M132").Value = -8488.750100000001
$ws.Range("N132").Value = -11876.6432

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 2730.8235
$ws.Range("J46").Value = 2994.9333
$ws.Range("L46").Value = 8984.7999
$ws.Range("N46").Value = -9166.7999

$ws.Range("H51").Value = 2813.7917
$ws.Range("I51").Value = 1500.3334
$ws.Range("K51").Value = 4501.0002
$ws.Range("M51").Value = -4041.0002

$ws.Range("H58").Value = 5000
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()

$ws.Range("H68").Value = 2629.8142
$ws.Range("I68").Value = 791.86664
$ws.Range("J68").Value = 4008.275
$ws.Range("K68").Value = 2375.59992
$ws.Range("L68").Value = 12024.825
$ws.Range("M68").Value = -1564.59992
$ws.Range("N68").Value = -13646.825

$ws.Range("H71").Value = 2629.8142
$ws.Range("I71").Value = 791.86664
$ws.Range("J71").Value = 4008.275
$ws.Range("K71").Value = 7126.79976
$ws.Range("L71").Value = 36074.475
$ws.Range("M71").Value = -3070.79976
$ws.Range("N71").Value = -44186.475

$ws.Range("H76").Value = 3270
$ws.Range("I76").Value = 2640
$ws.Range("J76").Value = 3900
$ws.Range("K76").Value = 7920
$ws.Range("L76").Value = 11700
$ws.Range("M76").Value = -7537
$ws.Range("N76").Value = -12466

$ws.Range("H79").Value = 3270
$ws.Range("I79").Value = 2640
$ws.Range("J79").Value = 3900
$ws.Range("K79").Value = 7920
$ws.Range("L79").Value = 11700
$ws.Range("M79").Value = -6594
$ws.Range("N79").Value = -14352

$ws.Range("H107").Value = 4457.4
$ws.Range("J107").Value = 2416.0952
$ws.Range("L107").Value = 7248.285600000001
$ws.Range("N107").Value = -11088.2856

$ws.Range("H132").Value = 9183.799999999999
$ws.Range("I132").Value = 4599.3335
$ws.Range("J132").Value = 10329.917
$ws.Range("K132").Value = 41394.0015
$ws.Range("L132").Value = 92969.253
$ws.Range("M132").Value = -38864.0015
$ws.Range("N132").Value = -98029.253

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1639.25
$ws.Range("I97").Value = 1696.3636
$ws.Range("J97").Value = 1011
$ws.Range("K97").Value = 1696.3636
$ws.Range("L97").Value = 1011
$ws.Range("M97").Value = -1200.3636
$ws.Range("N97").Value = -2003
